# Update Case_0_44 vm_pu.xlsx results: slack bus voltage setpoint changed from 1.05 to 1.02 p.u.
# (bus B is the slack/external grid bus), with recomputed per-unit voltages for all buses/rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.021515063812823
$ws.Cells.Item(2, 4).Value = 1.023273863702165
$ws.Cells.Item(2, 5).Value = 1.022383319796544
$ws.Cells.Item(2, 6).Value = 1.019947287561612
$ws.Cells.Item(2, 9).Value = 1.028946112253071
$ws.Cells.Item(2, 10).Value = 1.026706097584582
$ws.Cells.Item(2, 11).Value = 1.026105784809271
$ws.Cells.Item(2, 12).Value = 1.025217862320743
$ws.Cells.Item(2, 13).Value = 1.022789026279868
$ws.Cells.Item(2, 14).Value = 1.028164138007989

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022549113817173
$ws.Cells.Item(3, 4).Value = 1.024168149069848
$ws.Cells.Item(3, 5).Value = 1.023261920115625
$ws.Cells.Item(3, 6).Value = 1.021621363779628
$ws.Cells.Item(3, 9).Value = 1.029103721808994
$ws.Cells.Item(3, 10).Value = 1.027377359477933
$ws.Cells.Item(3, 11).Value = 1.026806675667701
$ws.Cells.Item(3, 12).Value = 1.025902920499346
$ws.Cells.Item(3, 13).Value = 1.02426685457582
$ws.Cells.Item(3, 14).Value = 1.028836353170221

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.023217414583669
$ws.Cells.Item(4, 4).Value = 1.024746389729723
$ws.Cells.Item(4, 5).Value = 1.02383014941081
$ws.Cells.Item(4, 6).Value = 1.022703640087793
$ws.Cells.Item(4, 9).Value = 1.029203216583911
$ws.Cells.Item(4, 10).Value = 1.027810370383511
$ws.Cells.Item(4, 11).Value = 1.027259150795972
$ws.Cells.Item(4, 12).Value = 1.026345285075286
$ws.Cells.Item(4, 13).Value = 1.025221701594134
$ws.Cells.Item(4, 14).Value = 1.029269979000952

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02349817832584
$ws.Cells.Item(5, 4).Value = 1.024989381707334
$ws.Cells.Item(5, 5).Value = 1.02406896551708
$ws.Cells.Item(5, 6).Value = 1.023158406189097
$ws.Cells.Item(5, 9).Value = 1.029244448131391
$ws.Cells.Item(5, 10).Value = 1.027992088016269
$ws.Cells.Item(5, 11).Value = 1.02744912107265
$ws.Cells.Item(5, 12).Value = 1.026531036937478
$ws.Cells.Item(5, 13).Value = 1.025622789140283
$ws.Cells.Item(5, 14).Value = 1.029451954693593

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.023545308695645
$ws.Cells.Item(6, 4).Value = 1.025030175258717
$ws.Cells.Item(6, 5).Value = 1.024109059847767
$ws.Cells.Item(6, 6).Value = 1.023234750572201
$ws.Cells.Item(6, 9).Value = 1.029251336131975
$ws.Cells.Item(6, 10).Value = 1.028022580453903
$ws.Cells.Item(6, 11).Value = 1.027481003260905
$ws.Cells.Item(6, 12).Value = 1.026562212713245
$ws.Cells.Item(6, 13).Value = 1.025690114336436
$ws.Cells.Item(6, 14).Value = 1.029482490433985

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.023221166907109
$ws.Cells.Item(7, 4).Value = 1.024749636994233
$ws.Cells.Item(7, 5).Value = 1.023833340747959
$ws.Cells.Item(7, 6).Value = 1.022709717560907
$ws.Cells.Item(7, 9).Value = 1.029203769864184
$ws.Cells.Item(7, 10).Value = 1.027812799760506
$ws.Cells.Item(7, 11).Value = 1.027261690170713
$ws.Cells.Item(7, 12).Value = 1.026347767957482
$ws.Cells.Item(7, 13).Value = 1.025227062234613
$ws.Cells.Item(7, 14).Value = 1.029272411827942

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021864692035658
$ws.Cells.Item(8, 4).Value = 1.023576179314074
$ws.Cells.Item(8, 5).Value = 1.022680305951224
$ws.Cells.Item(8, 6).Value = 1.020513253472642
$ws.Cells.Item(8, 9).Value = 1.028999892450378
$ws.Cells.Item(8, 10).Value = 1.026933231685321
$ws.Cells.Item(8, 11).Value = 1.026342871578539
$ws.Cells.Item(8, 12).Value = 1.025449570517039
$ws.Cells.Item(8, 13).Value = 1.023288760541173
$ws.Cells.Item(8, 14).Value = 1.028391594665202

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.019468232199361
$ws.Cells.Item(9, 4).Value = 1.021505139290784
$ws.Cells.Item(9, 5).Value = 1.020646314843734
$ws.Cells.Item(9, 6).Value = 1.016635083205301
$ws.Cells.Item(9, 9).Value = 1.028621575927863
$ws.Cells.Item(9, 10).Value = 1.025373019356986
$ws.Cells.Item(9, 11).Value = 1.024715737681574
$ws.Cells.Item(9, 12).Value = 1.023859805470929
$ws.Cells.Item(9, 13).Value = 1.019862153449528
$ws.Cells.Item(9, 14).Value = 1.026829166656403

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.01786635002673
$ws.Cells.Item(10, 4).Value = 1.020122213196972
$ws.Cells.Item(10, 5).Value = 1.019288808521223
$ws.Cells.Item(10, 6).Value = 1.01404396134186
$ws.Cells.Item(10, 9).Value = 1.028356555664771
$ws.Cells.Item(10, 10).Value = 1.024325899479398
$ws.Cells.Item(10, 11).Value = 1.023625517246087
$ws.Cells.Item(10, 12).Value = 1.022795196371548
$ws.Cells.Item(10, 13).Value = 1.017569885180381
$ws.Cells.Item(10, 14).Value = 1.025780559748482

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017171690162111
$ws.Cells.Item(11, 4).Value = 1.019522850633648
$ws.Cells.Item(11, 5).Value = 1.018700625476876
$ws.Cells.Item(11, 6).Value = 1.012920522416004
$ws.Cells.Item(11, 9).Value = 1.028238760882324
$ws.Cells.Item(11, 10).Value = 1.023870817047414
$ws.Cells.Item(11, 11).Value = 1.02315213215941
$ws.Cells.Item(11, 12).Value = 1.022333068180679
$ws.Cells.Item(11, 13).Value = 1.016575348495206
$ws.Cells.Item(11, 14).Value = 1.025324831047245

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016913505220619
$ws.Cells.Item(12, 4).Value = 1.019300137541749
$ws.Cells.Item(12, 5).Value = 1.018482091026367
$ws.Cells.Item(12, 6).Value = 1.0125029974856
$ws.Cells.Item(12, 9).Value = 1.028194549937481
$ws.Cells.Item(12, 10).Value = 1.023701526652524
$ws.Cells.Item(12, 11).Value = 1.022976097570714
$ws.Cells.Item(12, 12).Value = 1.022161240102123
$ws.Cells.Item(12, 13).Value = 1.016205629207559
$ws.Cells.Item(12, 14).Value = 1.025155300240577

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.016968893918378
$ws.Cells.Item(13, 4).Value = 1.019347914035834
$ws.Cells.Item(13, 5).Value = 1.018528970005685
$ws.Cells.Item(13, 6).Value = 1.012592568546764
$ws.Cells.Item(13, 9).Value = 1.028204054004515
$ws.Cells.Item(13, 10).Value = 1.023737851463767
$ws.Cells.Item(13, 11).Value = 1.02301386657867
$ws.Cells.Item(13, 12).Value = 1.022198105657418
$ws.Cells.Item(13, 13).Value = 1.01628494914246
$ws.Cells.Item(13, 14).Value = 1.025191676637218

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017150351735295
$ws.Cells.Item(14, 4).Value = 1.019504442802398
$ws.Cells.Item(14, 5).Value = 1.018682562520783
$ws.Cells.Item(14, 6).Value = 1.012886014415482
$ws.Cells.Item(14, 9).Value = 1.028235115709398
$ws.Cells.Item(14, 10).Value = 1.023856828609011
$ws.Cells.Item(14, 11).Value = 1.023137585136958
$ws.Cells.Item(14, 12).Value = 1.022318868355026
$ws.Cells.Item(14, 13).Value = 1.016544793634385
$ws.Cells.Item(14, 14).Value = 1.025310822743656

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.017262132977724
$ws.Cells.Item(15, 4).Value = 1.019600874291557
$ws.Cells.Item(15, 5).Value = 1.018777188347295
$ws.Cells.Item(15, 6).Value = 1.013066785380264
$ws.Cells.Item(15, 9).Value = 1.028254193328002
$ws.Cells.Item(15, 10).Value = 1.023930100856881
$ws.Cells.Item(15, 11).Value = 1.02321378591497
$ws.Cells.Item(15, 12).Value = 1.022393251271845
$ws.Cells.Item(15, 13).Value = 1.016704851861683
$ws.Cells.Item(15, 14).Value = 1.025384199046522

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.017912430357017
$ws.Cells.Item(16, 4).Value = 1.020161979314212
$ws.Cells.Item(16, 5).Value = 1.019327836338853
$ws.Cells.Item(16, 6).Value = 1.014118488591654
$ws.Cells.Item(16, 9).Value = 1.028364309274398
$ws.Cells.Item(16, 10).Value = 1.024356066457508
$ws.Cells.Item(16, 11).Value = 1.023656906509335
$ws.Cells.Item(16, 12).Value = 1.02282584207953
$ws.Cells.Item(16, 13).Value = 1.01763584711323
$ws.Cells.Item(16, 14).Value = 1.02581076956716

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.018320066442804
$ws.Cells.Item(17, 4).Value = 1.020513798637682
$ws.Cells.Item(17, 5).Value = 1.019673142549138
$ws.Cells.Item(17, 6).Value = 1.014777795086813
$ws.Cells.Item(17, 9).Value = 1.028432568334539
$ws.Cells.Item(17, 10).Value = 1.024622814760533
$ws.Cells.Item(17, 11).Value = 1.023934512195959
$ws.Cells.Item(17, 12).Value = 1.023096887498592
$ws.Cells.Item(17, 13).Value = 1.018219302803868
$ws.Cells.Item(17, 14).Value = 1.026077896683373

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.018557734012834
$ws.Cells.Item(18, 4).Value = 1.020718956247082
$ws.Cells.Item(18, 5).Value = 1.019874517812693
$ws.Cells.Item(18, 6).Value = 1.01516221659883
$ws.Cells.Item(18, 9).Value = 1.028472089453626
$ws.Cells.Item(18, 10).Value = 1.02477824323471
$ws.Cells.Item(18, 11).Value = 1.024096308341778
$ws.Cells.Item(18, 12).Value = 1.023254873258976
$ws.Cells.Item(18, 13).Value = 1.01855943317656
$ws.Cells.Item(18, 14).Value = 1.026233545883811

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.018638755665812
$ws.Cells.Item(19, 4).Value = 1.020788900790837
$ws.Cells.Item(19, 5).Value = 1.019943175492955
$ws.Cells.Item(19, 6).Value = 1.015293270799914
$ws.Cells.Item(19, 9).Value = 1.028485515398635
$ws.Cells.Item(19, 10).Value = 1.024831213005151
$ws.Cells.Item(19, 11).Value = 1.024151455213306
$ws.Cells.Item(19, 12).Value = 1.023308723616522
$ws.Cells.Item(19, 13).Value = 1.018675376954946
$ws.Cells.Item(19, 14).Value = 1.026286590877401

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.018276341272244
$ws.Cells.Item(20, 4).Value = 1.020476057208096
$ws.Cells.Item(20, 5).Value = 1.019636098192455
$ws.Cells.Item(20, 6).Value = 1.014707072374697
$ws.Cells.Item(20, 9).Value = 1.028425275115223
$ws.Cells.Item(20, 10).Value = 1.024594211885503
$ws.Cells.Item(20, 11).Value = 1.023904740831896
$ws.Cells.Item(20, 12).Value = 1.023067818322271
$ws.Cells.Item(20, 13).Value = 1.018156723188046
$ws.Cells.Item(20, 14).Value = 1.026049253188979

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017096921289978
$ws.Cells.Item(21, 4).Value = 1.019458351285932
$ws.Cells.Item(21, 5).Value = 1.018637334936388
$ws.Cells.Item(21, 6).Value = 1.012799608337833
$ws.Cells.Item(21, 9).Value = 1.028225981420068
$ws.Cells.Item(21, 10).Value = 1.023821799785838
$ws.Cells.Item(21, 11).Value = 1.023101158587955
$ws.Cells.Item(21, 12).Value = 1.022283311536887
$ws.Cells.Item(21, 13).Value = 1.016468284368272
$ws.Cells.Item(21, 14).Value = 1.025275744175536

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.016354461844577
$ws.Cells.Item(22, 4).Value = 1.018817997564172
$ws.Cells.Item(22, 5).Value = 1.018009042125687
$ws.Cells.Item(22, 6).Value = 1.011598976339489
$ws.Cells.Item(22, 9).Value = 1.028098035071606
$ws.Cells.Item(22, 10).Value = 1.023334692030949
$ws.Cells.Item(22, 11).Value = 1.022594766725152
$ws.Cells.Item(22, 12).Value = 1.021789059011813
$ws.Cells.Item(22, 13).Value = 1.015404932441854
$ws.Cells.Item(22, 14).Value = 1.024787944671762

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.016748140680841
$ws.Cells.Item(23, 4).Value = 1.019157507198113
$ws.Cells.Item(23, 5).Value = 1.018342143691152
$ws.Cells.Item(23, 6).Value = 1.012235583663137
$ws.Cells.Item(23, 9).Value = 1.028166112368706
$ws.Cells.Item(23, 10).Value = 1.023593056027996
$ws.Cells.Item(23, 11).Value = 1.022863323809365
$ws.Cells.Item(23, 12).Value = 1.022051166944006
$ws.Cells.Item(23, 13).Value = 1.015968805326168
$ws.Cells.Item(23, 14).Value = 1.025046675575321

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.01829609910266
$ws.Cells.Item(24, 4).Value = 1.020493111098335
$ws.Cells.Item(24, 5).Value = 1.019652837054043
$ws.Cells.Item(24, 6).Value = 1.014739029358384
$ws.Cells.Item(24, 9).Value = 1.028428571513214
$ws.Cells.Item(24, 10).Value = 1.024607136791946
$ws.Cells.Item(24, 11).Value = 1.023918193620901
$ws.Cells.Item(24, 12).Value = 1.023080953773099
$ws.Cells.Item(24, 13).Value = 1.018185000807461
$ws.Cells.Item(24, 14).Value = 1.026062196450272

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.020088515645946
$ws.Cells.Item(25, 4).Value = 1.022040943110268
$ws.Cells.Item(25, 5).Value = 1.021172414156683
$ws.Cells.Item(25, 6).Value = 1.017638647706223
$ws.Cells.Item(25, 9).Value = 1.028721637802818
$ws.Cells.Item(25, 10).Value = 1.025777597680652
$ws.Cells.Item(25, 11).Value = 1.025137350445669
$ws.Cells.Item(25, 12).Value = 1.024271634705388
$ws.Cells.Item(25, 13).Value = 1.020749368581536
$ws.Cells.Item(25, 14).Value = 1.027234319527694
